# Update sets and tables
# Insert a new "medium goods truck" demand-set row into the VEDA_Sets-Proc
# sheet, between the existing "light goods truck" (row 5) and "train"
# (formerly row 6, now row 7) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("VEDA_Sets-Proc")
$ws.Activate()

# Push rows 6:48 down to 7:49 by inserting a blank row at row 6.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 following the same pattern as the other
# DMD_TRA-F-* rows above it (heavy / land / light goods truck).
$ws.Cells.Item(6, 1).Value = "DMD"
$ws.Cells.Item(6, 2).Value = "T-MGT*"
$ws.Cells.Item(6, 5).Value = "TRAF*"
$ws.Cells.Item(6, 6).Value = "DMD_TRA-F-MTRUCK"
$ws.Cells.Item(6, 7).Value = "Freight transport - medium goods truck"
$ws.Cells.Item(6, 8).Value = "AND"
$ws.Cells.Item(6, 9).Value = "AND"
$ws.Cells.Item(6, 10).Value = "AND"
$ws.Cells.Item(6, 11).Value = "OR"

# Restore the view: no frozen/scrolled top-left cell and selection on G5.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G5").Select()
